# Update the Sterling Equity performance table (Sheet1) with the revised
# figures for the "Year to 30 June 2008" and "Year to 30 June 2009" rows,
# as requested.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 3 (Year to 30 June 2008): Sterling Equity / Index / Relative performance
$ws.Range("C3").Value = -0.059
$ws.Range("D3").Value = -0.137
$ws.Range("E3").Value = 0.078

# Row 4 (Year to 30 June 2009): Sterling Equity / Index / Relative performance
$ws.Range("C4").Value = 0.361
$ws.Range("D4").Value = -0.203
$ws.Range("E4").Value = 0.565

# Bring the view back to the top of the sheet and leave the selection on
# the bottom annualised-performance row, matching the refreshed table view.
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("E11").Select()
